# Append 45 new data rows (rows 102-146) to the
# "master-reg_center_device_h" worksheet, continuing the existing
# regcntr_id / device_id cycling pattern, and refresh the sheet's
# selection / print-orientation view state to match the edited file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id (column A) cycles through this 9-value sequence.
$cycle = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)

$startRow = 102
$endRow = 146
$startDeviceId = 3000121

for ($i = 0; $i -le ($endRow - $startRow); $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $cycle[$i % 9]
    $ws.Cells.Item($row, 2).Value = $startDeviceId + $i
}

# Columns C-G are identical on every new row, so broadcast them across
# the whole new block in one shot each instead of cell-by-cell.
$ws.Range("C$startRow`:C$endRow").Value = "eng"
$ws.Range("D$startRow`:D$endRow").Value = $true
$ws.Range("E$startRow`:E$endRow").Value = "superadmin"
$ws.Range("F$startRow`:G$endRow").Value = "now()"

# Match the saved file's selection / scroll state: the new block is
# selected with the active cell at its top-left corner.
$ws.Range("A102:B146").Select()

# Match the saved file's print setup (portrait orientation).
$ws.PageSetup.Orientation = 1
